# Append: 2025-10-31 01:53 JST
# Update the "取得日時" (acquisition datetime) column for all existing
# data rows on the "ランサーズ" sheet to reflect the latest scrape run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-31 01:53:38"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
